$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("dSF") value corrections per repulled data
$ws.Range("F4").Value = 10
$ws.Range("F6").Value = -15
$ws.Range("F7").Value = -5
$ws.Range("F8").Value = -12
$ws.Range("F9").Value = -6
$ws.Range("F10").Value = -3
$ws.Range("F11").Value = -5
$ws.Range("F13").Value = -10
$ws.Range("F16").Value = -5
$ws.Range("F18").Value = 6
$ws.Range("F19").Value = -3
$ws.Range("F24").Value = -2
$ws.Range("F26").Value = -5
$ws.Range("F27").Value = -7
$ws.Range("F33").Value = 1
